$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item(1)

# --- Step 1: Insert new row 2 into the "总计" (summary) sheet for 2022-Q3 ---
$zj.Rows.Item(2).Insert()
$zj.Rows.Item(2).ClearFormats()
$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q3"
$zj.Cells.Item(2,3).Value = 46
$zj.Cells.Item(2,4).Value = 12.57
# Restore the A2 style (bold/centered/bordered), copying it from A3 which still carries the original style
$zj.Cells.Item(3,1).Copy()
$zj.Cells.Item(2,1).PasteSpecial(-4122)

# --- Step 2: Create the new "2022-Q3" worksheet positioned right after "总计" ---
$q3 = $wb.Worksheets.Add($null, $zj)
$q3.Name = "2022-Q3"

# Reference sheet used only to copy matching header/index cell styles (layout identical across quarter sheets)
$ref = $wb.Worksheets.Item("2022-Q2")

# Header row values
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Row 2: 011164 富国兴远优选12个月持有期混合A
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "'011164"
$q3.Cells.Item(2,3).Value = "富国兴远优选12个月持有期混合A"
$q3.Cells.Item(2,4).Value = "'37.00"
$q3.Cells.Item(2,5).Value = "'73.75"
$q3.Cells.Item(2,6).Value = "'9.76"
$q3.Cells.Item(2,7).Value = "'3.6112"
$q3.Cells.Item(2,8).Value = 1
# Row 3: 011165 富国兴远优选12个月持有期混合C
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "'011165"
$q3.Cells.Item(3,3).Value = "富国兴远优选12个月持有期混合C"
$q3.Cells.Item(3,4).Value = "'15.92"
$q3.Cells.Item(3,5).Value = "'73.75"
$q3.Cells.Item(3,6).Value = "'9.76"
$q3.Cells.Item(3,7).Value = "'1.5538"
$q3.Cells.Item(3,8).Value = 1
# Row 4: 110002 易方达策略成长混合
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "'110002"
$q3.Cells.Item(4,3).Value = "易方达策略成长混合"
$q3.Cells.Item(4,4).Value = "'11.73"
$q3.Cells.Item(4,5).Value = "'90.54"
$q3.Cells.Item(4,6).Value = "'6.73"
$q3.Cells.Item(4,7).Value = "'0.7894"
$q3.Cells.Item(4,8).Value = 2
# Row 5: 001186 富国文体健康股票A
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "'001186"
$q3.Cells.Item(5,3).Value = "富国文体健康股票A"
$q3.Cells.Item(5,4).Value = "'9.91"
$q3.Cells.Item(5,5).Value = "'82.41"
$q3.Cells.Item(5,6).Value = "'7.71"
$q3.Cells.Item(5,7).Value = "'0.7641"
$q3.Cells.Item(5,8).Value = 3
# Row 6: 010846 南方卓越优选3个月持有期混合A
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "'010846"
$q3.Cells.Item(6,3).Value = "南方卓越优选3个月持有期混合A"
$q3.Cells.Item(6,4).Value = "'20.44"
$q3.Cells.Item(6,5).Value = "'80.88"
$q3.Cells.Item(6,6).Value = "'3.54"
$q3.Cells.Item(6,7).Value = "'0.7236"
$q3.Cells.Item(6,8).Value = 10
# Row 7: 112002 易方达策略成长二号混合
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "'112002"
$q3.Cells.Item(7,3).Value = "易方达策略成长二号混合"
$q3.Cells.Item(7,4).Value = "'9.94"
$q3.Cells.Item(7,5).Value = "'91.00"
$q3.Cells.Item(7,6).Value = "'6.74"
$q3.Cells.Item(7,7).Value = "'0.6700"
$q3.Cells.Item(7,8).Value = 2
# Row 8: 202023 南方优选成长混合A
$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).Value = "'202023"
$q3.Cells.Item(8,3).Value = "南方优选成长混合A"
$q3.Cells.Item(8,4).Value = "'37.42"
$q3.Cells.Item(8,5).Value = "'73.42"
$q3.Cells.Item(8,6).Value = "'1.79"
$q3.Cells.Item(8,7).Value = "'0.6698"
$q3.Cells.Item(8,8).Value = 8
# Row 9: 005123 南方优享分红灵活配置混合A
$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).Value = "'005123"
$q3.Cells.Item(9,3).Value = "南方优享分红灵活配置混合A"
$q3.Cells.Item(9,4).Value = "'7.37"
$q3.Cells.Item(9,5).Value = "'92.25"
$q3.Cells.Item(9,6).Value = "'7.79"
$q3.Cells.Item(9,7).Value = "'0.5741"
$q3.Cells.Item(9,8).Value = 4
# Row 10: 008854 南方内需增长两年持有期股票A
$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).Value = "'008854"
$q3.Cells.Item(10,3).Value = "南方内需增长两年持有期股票A"
$q3.Cells.Item(10,4).Value = "'21.14"
$q3.Cells.Item(10,5).Value = "'88.74"
$q3.Cells.Item(10,6).Value = "'2.05"
$q3.Cells.Item(10,7).Value = "'0.4334"
$q3.Cells.Item(10,8).Value = 10
# Row 11: 515210 国泰中证钢铁ETF
$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).Value = "'515210"
$q3.Cells.Item(11,3).Value = "国泰中证钢铁ETF"
$q3.Cells.Item(11,4).Value = "'14.23"
$q3.Cells.Item(11,5).Value = "'97.88"
$q3.Cells.Item(11,6).Value = "'2.93"
$q3.Cells.Item(11,7).Value = "'0.4169"
$q3.Cells.Item(11,8).Value = 7
# Row 12: 002851 南方品质优选灵活配置混合A
$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).Value = "'002851"
$q3.Cells.Item(12,3).Value = "南方品质优选灵活配置混合A"
$q3.Cells.Item(12,4).Value = "'11.81"
$q3.Cells.Item(12,5).Value = "'71.72"
$q3.Cells.Item(12,6).Value = "'3.30"
$q3.Cells.Item(12,7).Value = "'0.3897"
$q3.Cells.Item(12,8).Value = 9
# Row 13: 502023 鹏华国证钢铁行业指数（LOF）A
$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).Value = "'502023"
$q3.Cells.Item(13,3).Value = "鹏华国证钢铁行业指数（LOF）A"
$q3.Cells.Item(13,4).Value = "'9.48"
$q3.Cells.Item(13,5).Value = "'94.49"
$q3.Cells.Item(13,6).Value = "'2.90"
$q3.Cells.Item(13,7).Value = "'0.2749"
$q3.Cells.Item(13,8).Value = 7
# Row 14: 011830 富国天恒混合A
$q3.Cells.Item(14,1).Value = 12
$q3.Cells.Item(14,2).Value = "'011830"
$q3.Cells.Item(14,3).Value = "富国天恒混合A"
$q3.Cells.Item(14,4).Value = "'3.26"
$q3.Cells.Item(14,5).Value = "'65.79"
$q3.Cells.Item(14,6).Value = "'8.29"
$q3.Cells.Item(14,7).Value = "'0.2703"
$q3.Cells.Item(14,8).Value = 1
# Row 15: 012810 鹏华国证钢铁行业指数（LOF）C
$q3.Cells.Item(15,1).Value = 13
$q3.Cells.Item(15,2).Value = "'012810"
$q3.Cells.Item(15,3).Value = "鹏华国证钢铁行业指数（LOF）C"
$q3.Cells.Item(15,4).Value = "'4.34"
$q3.Cells.Item(15,5).Value = "'94.49"
$q3.Cells.Item(15,6).Value = "'2.90"
$q3.Cells.Item(15,7).Value = "'0.1259"
$q3.Cells.Item(15,8).Value = 7
# Row 16: 010847 南方卓越优选3个月持有期混合C
$q3.Cells.Item(16,1).Value = 14
$q3.Cells.Item(16,2).Value = "'010847"
$q3.Cells.Item(16,3).Value = "南方卓越优选3个月持有期混合C"
$q3.Cells.Item(16,4).Value = "'3.06"
$q3.Cells.Item(16,5).Value = "'80.88"
$q3.Cells.Item(16,6).Value = "'3.54"
$q3.Cells.Item(16,7).Value = "'0.1083"
$q3.Cells.Item(16,8).Value = 10
# Row 17: 008855 南方内需增长两年持有期股票C
$q3.Cells.Item(17,1).Value = 15
$q3.Cells.Item(17,2).Value = "'008855"
$q3.Cells.Item(17,3).Value = "南方内需增长两年持有期股票C"
$q3.Cells.Item(17,4).Value = "'5.07"
$q3.Cells.Item(17,5).Value = "'88.74"
$q3.Cells.Item(17,6).Value = "'2.05"
$q3.Cells.Item(17,7).Value = "'0.1039"
$q3.Cells.Item(17,8).Value = 10
# Row 18: 168203 中融国证钢铁行业指数A
$q3.Cells.Item(18,1).Value = 16
$q3.Cells.Item(18,2).Value = "'168203"
$q3.Cells.Item(18,3).Value = "中融国证钢铁行业指数A"
$q3.Cells.Item(18,4).Value = "'3.34"
$q3.Cells.Item(18,5).Value = "'92.81"
$q3.Cells.Item(18,6).Value = "'2.85"
$q3.Cells.Item(18,7).Value = "'0.0952"
$q3.Cells.Item(18,8).Value = 7
# Row 19: 001443 易方达瑞选灵活配置混合I
$q3.Cells.Item(19,1).Value = 17
$q3.Cells.Item(19,2).Value = "'001443"
$q3.Cells.Item(19,3).Value = "易方达瑞选灵活配置混合I"
$q3.Cells.Item(19,4).Value = "'5.84"
$q3.Cells.Item(19,5).Value = "'28.98"
$q3.Cells.Item(19,6).Value = "'1.49"
$q3.Cells.Item(19,7).Value = "'0.0870"
$q3.Cells.Item(19,8).Value = 5
# Row 20: 001444 易方达瑞选灵活配置混合E
$q3.Cells.Item(20,1).Value = 18
$q3.Cells.Item(20,2).Value = "'001444"
$q3.Cells.Item(20,3).Value = "易方达瑞选灵活配置混合E"
$q3.Cells.Item(20,4).Value = "'5.84"
$q3.Cells.Item(20,5).Value = "'28.98"
$q3.Cells.Item(20,6).Value = "'1.49"
$q3.Cells.Item(20,7).Value = "'0.0870"
$q3.Cells.Item(20,8).Value = 5
# Row 21: 011142 创金合信新材料新能源股票A
$q3.Cells.Item(21,1).Value = 19
$q3.Cells.Item(21,2).Value = "'011142"
$q3.Cells.Item(21,3).Value = "创金合信新材料新能源股票A"
$q3.Cells.Item(21,4).Value = "'2.04"
$q3.Cells.Item(21,5).Value = "'89.16"
$q3.Cells.Item(21,6).Value = "'4.12"
$q3.Cells.Item(21,7).Value = "'0.0840"
$q3.Cells.Item(21,8).Value = 2
# Row 22: 014031 南方发展机遇一年持有混合A
$q3.Cells.Item(22,1).Value = 20
$q3.Cells.Item(22,2).Value = "'014031"
$q3.Cells.Item(22,3).Value = "南方发展机遇一年持有混合A"
$q3.Cells.Item(22,4).Value = "'3.48"
$q3.Cells.Item(22,5).Value = "'83.48"
$q3.Cells.Item(22,6).Value = "'2.11"
$q3.Cells.Item(22,7).Value = "'0.0734"
$q3.Cells.Item(22,8).Value = 7
# Row 23: 009215 易方达瑞川灵活配置混合A
$q3.Cells.Item(23,1).Value = 21
$q3.Cells.Item(23,2).Value = "'009215"
$q3.Cells.Item(23,3).Value = "易方达瑞川灵活配置混合A"
$q3.Cells.Item(23,4).Value = "'6.45"
$q3.Cells.Item(23,5).Value = "'28.88"
$q3.Cells.Item(23,6).Value = "'1.11"
$q3.Cells.Item(23,7).Value = "'0.0716"
$q3.Cells.Item(23,8).Value = 9
# Row 24: 001314 易方达新益灵活配置混合I
$q3.Cells.Item(24,1).Value = 22
$q3.Cells.Item(24,2).Value = "'001314"
$q3.Cells.Item(24,3).Value = "易方达新益灵活配置混合I"
$q3.Cells.Item(24,4).Value = "'7.02"
$q3.Cells.Item(24,5).Value = "'24.50"
$q3.Cells.Item(24,6).Value = "'1.01"
$q3.Cells.Item(24,7).Value = "'0.0709"
$q3.Cells.Item(24,8).Value = 6
# Row 25: 001315 易方达新益灵活配置混合E
$q3.Cells.Item(25,1).Value = 23
$q3.Cells.Item(25,2).Value = "'001315"
$q3.Cells.Item(25,3).Value = "易方达新益灵活配置混合E"
$q3.Cells.Item(25,4).Value = "'7.02"
$q3.Cells.Item(25,5).Value = "'24.50"
$q3.Cells.Item(25,6).Value = "'1.01"
$q3.Cells.Item(25,7).Value = "'0.0709"
$q3.Cells.Item(25,8).Value = 6
# Row 26: 001747 易方达瑞祺灵活配置混合I
$q3.Cells.Item(26,1).Value = 24
$q3.Cells.Item(26,2).Value = "'001747"
$q3.Cells.Item(26,3).Value = "易方达瑞祺灵活配置混合I"
$q3.Cells.Item(26,4).Value = "'5.14"
$q3.Cells.Item(26,5).Value = "'28.39"
$q3.Cells.Item(26,6).Value = "'1.33"
$q3.Cells.Item(26,7).Value = "'0.0684"
$q3.Cells.Item(26,8).Value = 5
# Row 27: 001748 易方达瑞祺灵活配置混合E
$q3.Cells.Item(27,1).Value = 25
$q3.Cells.Item(27,2).Value = "'001748"
$q3.Cells.Item(27,3).Value = "易方达瑞祺灵活配置混合E"
$q3.Cells.Item(27,4).Value = "'5.14"
$q3.Cells.Item(27,5).Value = "'28.39"
$q3.Cells.Item(27,6).Value = "'1.33"
$q3.Cells.Item(27,7).Value = "'0.0684"
$q3.Cells.Item(27,8).Value = 5
# Row 28: 011903 南方领航优选混合A
$q3.Cells.Item(28,1).Value = 26
$q3.Cells.Item(28,2).Value = "'011903"
$q3.Cells.Item(28,3).Value = "南方领航优选混合A"
$q3.Cells.Item(28,4).Value = "'1.56"
$q3.Cells.Item(28,5).Value = "'82.01"
$q3.Cells.Item(28,6).Value = "'3.99"
$q3.Cells.Item(28,7).Value = "'0.0622"
$q3.Cells.Item(28,8).Value = 9
# Row 29: 011143 创金合信新材料新能源股票C
$q3.Cells.Item(29,1).Value = 27
$q3.Cells.Item(29,2).Value = "'011143"
$q3.Cells.Item(29,3).Value = "创金合信新材料新能源股票C"
$q3.Cells.Item(29,4).Value = "'1.19"
$q3.Cells.Item(29,5).Value = "'89.16"
$q3.Cells.Item(29,6).Value = "'4.12"
$q3.Cells.Item(29,7).Value = "'0.0490"
$q3.Cells.Item(29,8).Value = 2
# Row 30: 004703 南方兴盛先锋灵活配置混合
$q3.Cells.Item(30,1).Value = 28
$q3.Cells.Item(30,2).Value = "'004703"
$q3.Cells.Item(30,3).Value = "南方兴盛先锋灵活配置混合"
$q3.Cells.Item(30,4).Value = "'0.97"
$q3.Cells.Item(30,5).Value = "'82.88"
$q3.Cells.Item(30,6).Value = "'4.22"
$q3.Cells.Item(30,7).Value = "'0.0409"
$q3.Cells.Item(30,8).Value = 9
# Row 31: 006587 南方优享分红灵活配置混合C
$q3.Cells.Item(31,1).Value = 29
$q3.Cells.Item(31,2).Value = "'006587"
$q3.Cells.Item(31,3).Value = "南方优享分红灵活配置混合C"
$q3.Cells.Item(31,4).Value = "'0.46"
$q3.Cells.Item(31,5).Value = "'92.25"
$q3.Cells.Item(31,6).Value = "'7.79"
$q3.Cells.Item(31,7).Value = "'0.0358"
$q3.Cells.Item(31,8).Value = 4
# Row 32: 005206 南方优选成长混合C
$q3.Cells.Item(32,1).Value = 30
$q3.Cells.Item(32,2).Value = "'005206"
$q3.Cells.Item(32,3).Value = "南方优选成长混合C"
$q3.Cells.Item(32,4).Value = "'1.88"
$q3.Cells.Item(32,5).Value = "'73.42"
$q3.Cells.Item(32,6).Value = "'1.79"
$q3.Cells.Item(32,7).Value = "'0.0337"
$q3.Cells.Item(32,8).Value = 8
# Row 33: 014032 南方发展机遇一年持有混合C
$q3.Cells.Item(33,1).Value = 31
$q3.Cells.Item(33,2).Value = "'014032"
$q3.Cells.Item(33,3).Value = "南方发展机遇一年持有混合C"
$q3.Cells.Item(33,4).Value = "'1.14"
$q3.Cells.Item(33,5).Value = "'83.48"
$q3.Cells.Item(33,6).Value = "'2.11"
$q3.Cells.Item(33,7).Value = "'0.0241"
$q3.Cells.Item(33,8).Value = 7
# Row 34: 011904 南方领航优选混合C
$q3.Cells.Item(34,1).Value = 32
$q3.Cells.Item(34,2).Value = "'011904"
$q3.Cells.Item(34,3).Value = "南方领航优选混合C"
$q3.Cells.Item(34,4).Value = "'0.48"
$q3.Cells.Item(34,5).Value = "'82.01"
$q3.Cells.Item(34,6).Value = "'3.99"
$q3.Cells.Item(34,7).Value = "'0.0192"
$q3.Cells.Item(34,8).Value = 9
# Row 35: 011359 长城优选添利一年持有期混合A
$q3.Cells.Item(35,1).Value = 33
$q3.Cells.Item(35,2).Value = "'011359"
$q3.Cells.Item(35,3).Value = "长城优选添利一年持有期混合A"
$q3.Cells.Item(35,4).Value = "'0.99"
$q3.Cells.Item(35,5).Value = "'28.41"
$q3.Cells.Item(35,6).Value = "'1.27"
$q3.Cells.Item(35,7).Value = "'0.0126"
$q3.Cells.Item(35,8).Value = 5
# Row 36: 011125 富国文体健康股票C
$q3.Cells.Item(36,1).Value = 34
$q3.Cells.Item(36,2).Value = "'011125"
$q3.Cells.Item(36,3).Value = "富国文体健康股票C"
$q3.Cells.Item(36,4).Value = "'0.14"
$q3.Cells.Item(36,5).Value = "'82.41"
$q3.Cells.Item(36,6).Value = "'7.71"
$q3.Cells.Item(36,7).Value = "'0.0108"
$q3.Cells.Item(36,8).Value = 3
# Row 37: 006182 格林伯锐灵活配置混合C
$q3.Cells.Item(37,1).Value = 35
$q3.Cells.Item(37,2).Value = "'006182"
$q3.Cells.Item(37,3).Value = "格林伯锐灵活配置混合C"
$q3.Cells.Item(37,4).Value = "'0.14"
$q3.Cells.Item(37,5).Value = "'85.60"
$q3.Cells.Item(37,6).Value = "'4.85"
$q3.Cells.Item(37,7).Value = "'0.0068"
$q3.Cells.Item(37,8).Value = 3
# Row 38: 006181 格林伯锐灵活配置混合A
$q3.Cells.Item(38,1).Value = 36
$q3.Cells.Item(38,2).Value = "'006181"
$q3.Cells.Item(38,3).Value = "格林伯锐灵活配置混合A"
$q3.Cells.Item(38,4).Value = "'0.13"
$q3.Cells.Item(38,5).Value = "'85.60"
$q3.Cells.Item(38,6).Value = "'4.85"
$q3.Cells.Item(38,7).Value = "'0.0063"
$q3.Cells.Item(38,8).Value = 3
# Row 39: 011831 富国天恒混合C
$q3.Cells.Item(39,1).Value = 37
$q3.Cells.Item(39,2).Value = "'011831"
$q3.Cells.Item(39,3).Value = "富国天恒混合C"
$q3.Cells.Item(39,4).Value = "'0.04"
$q3.Cells.Item(39,5).Value = "'65.79"
$q3.Cells.Item(39,6).Value = "'8.29"
$q3.Cells.Item(39,7).Value = "'0.0033"
$q3.Cells.Item(39,8).Value = 1
# Row 40: 014692 中加量化研选混合型证券投资基金C
$q3.Cells.Item(40,1).Value = 38
$q3.Cells.Item(40,2).Value = "'014692"
$q3.Cells.Item(40,3).Value = "中加量化研选混合型证券投资基金C"
$q3.Cells.Item(40,4).Value = "'0.19"
$q3.Cells.Item(40,5).Value = "'68.41"
$q3.Cells.Item(40,6).Value = "'1.70"
$q3.Cells.Item(40,7).Value = "'0.0032"
$q3.Cells.Item(40,8).Value = 10
# Row 41: 013802 财通资管中证钢铁指数A
$q3.Cells.Item(41,1).Value = 39
$q3.Cells.Item(41,2).Value = "'013802"
$q3.Cells.Item(41,3).Value = "财通资管中证钢铁指数A"
$q3.Cells.Item(41,4).Value = "'0.08"
$q3.Cells.Item(41,5).Value = "'92.45"
$q3.Cells.Item(41,6).Value = "'2.96"
$q3.Cells.Item(41,7).Value = "'0.0024"
$q3.Cells.Item(41,8).Value = 8
# Row 42: 009216 易方达瑞川灵活配置混合C
$q3.Cells.Item(42,1).Value = 40
$q3.Cells.Item(42,2).Value = "'009216"
$q3.Cells.Item(42,3).Value = "易方达瑞川灵活配置混合C"
$q3.Cells.Item(42,4).Value = "'0.21"
$q3.Cells.Item(42,5).Value = "'28.88"
$q3.Cells.Item(42,6).Value = "'1.11"
$q3.Cells.Item(42,7).Value = "'0.0023"
$q3.Cells.Item(42,8).Value = 9
# Row 43: 011360 长城优选添利一年持有期混合C
$q3.Cells.Item(43,1).Value = 41
$q3.Cells.Item(43,2).Value = "'011360"
$q3.Cells.Item(43,3).Value = "长城优选添利一年持有期混合C"
$q3.Cells.Item(43,4).Value = "'0.16"
$q3.Cells.Item(43,5).Value = "'28.41"
$q3.Cells.Item(43,6).Value = "'1.27"
$q3.Cells.Item(43,7).Value = "'0.0020"
$q3.Cells.Item(43,8).Value = 5
# Row 44: 014691 中加量化研选混合型证券投资基金A
$q3.Cells.Item(44,1).Value = 42
$q3.Cells.Item(44,2).Value = "'014691"
$q3.Cells.Item(44,3).Value = "中加量化研选混合型证券投资基金A"
$q3.Cells.Item(44,4).Value = "'0.05"
$q3.Cells.Item(44,5).Value = "'68.41"
$q3.Cells.Item(44,6).Value = "'1.70"
$q3.Cells.Item(44,7).Value = "'0.0008"
$q3.Cells.Item(44,8).Value = 10
# Row 45: 013803 财通资管中证钢铁指数C
$q3.Cells.Item(45,1).Value = 43
$q3.Cells.Item(45,2).Value = "'013803"
$q3.Cells.Item(45,3).Value = "财通资管中证钢铁指数C"
$q3.Cells.Item(45,4).Value = "'0.01"
$q3.Cells.Item(45,5).Value = "'92.45"
$q3.Cells.Item(45,6).Value = "'2.96"
$q3.Cells.Item(45,7).Value = "'0.0003"
$q3.Cells.Item(45,8).Value = 8
# Row 46: 013501 南方品质优选灵活配置混合C
$q3.Cells.Item(46,1).Value = 44
$q3.Cells.Item(46,2).Value = "'013501"
$q3.Cells.Item(46,3).Value = "南方品质优选灵活配置混合C"
$q3.Cells.Item(46,4).Value = "'0.00"
$q3.Cells.Item(46,5).Value = "'71.72"
$q3.Cells.Item(46,6).Value = "'3.30"
$q3.Cells.Item(46,7).Value = 0
$q3.Cells.Item(46,8).Value = 9
# Row 47: 016815 中融国证钢铁行业指数C
$q3.Cells.Item(47,1).Value = 45
$q3.Cells.Item(47,2).Value = "'016815"
$q3.Cells.Item(47,3).Value = "中融国证钢铁行业指数C"
$q3.Cells.Item(47,4).Value = "'0.00"
$q3.Cells.Item(47,5).Value = "'92.81"
$q3.Cells.Item(47,6).Value = "'2.85"
$q3.Cells.Item(47,7).Value = 0
$q3.Cells.Item(47,8).Value = 7

# --- Step 3: Apply matching cell styles (format-only paste; does not touch the values set above) ---
$ref.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$ref.Range("A2").Copy()
$q3.Range("A2:A47").PasteSpecial(-4122)

Write-Host "done"